$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 0.5603230236742718),
    @("C2", 0.1735613149463688),
    @("D2", 0.605397865688289),
    @("E2", 0.2248957411965691),
    @("G2", 1.827928899758618),
    @("H2", 1.585733382079042),
    @("I2", 1.299312586894857),
    @("J2", 0.1010246523313292),
    @("K2", 0.7238254563052351),
    @("M2", 0.4283839459954351),
    @("N2", 2.677684805392005),
    @("B3", 0.5275487372977636),
    @("C3", 0.1661058282110872),
    @("D3", 0.6003746456235888),
    @("E3", 0.2236264274368729),
    @("G3", 1.823640456781405),
    @("H3", 1.589033429581221),
    @("I3", 1.301685517700349),
    @("J3", 0.1008446484672731),
    @("K3", 0.6846736824159905),
    @("M3", 0.4170279941517094),
    @("N3", 2.696603945128501),
    @("B4", 0.5077129616831826),
    @("C4", 0.1616269454579964),
    @("D4", 0.5975704282415251),
    @("E4", 0.2229466068047294),
    @("G4", 1.821889979805277),
    @("H4", 1.591661091706257),
    @("I4", 1.303680209156212),
    @("J4", 0.1007749079319673),
    @("K4", 0.6610188743377989),
    @("M4", 0.4102842353495149),
    @("N4", 2.70895380723092),
    @("B5", 0.4997022907425048),
    @("C5", 0.1598265915468318),
    @("D5", 0.5964981816095047),
    @("E5", 0.2226946274598944),
    @("G5", 1.821398471675678),
    @("H5", 1.592883187622647),
    @("I5", 1.304628259855946),
    @("J5", 0.1007567503565738),
    @("K5", 0.6514762080888659),
    @("M5", 0.4075937102281628),
    @("N5", 2.714170855463337),
    @("B6", 0.4983765140726177),
    @("C6", 0.1595291430648018),
    @("D6", 0.5963243957114486),
    @("E6", 0.2226543002882977),
    @("G6", 1.82133025133237),
    @("H6", 1.593095256121558),
    @("I6", 1.304793848260616),
    @("J6", 0.1007543553348818),
    @("K6", 0.649897510691261),
    @("M6", 0.4071504321850981),
    @("N6", 2.715048277024444),
    @("B7", 0.5076046328700556),
    @("C7", 0.1616025647584962),
    @("D7", 0.597555682021337),
    @("E7", 0.2229431070542915),
    @("G7", 1.821882453158665),
    @("H7", 1.591676960613754),
    @("I7", 1.303692447508567),
    @("J7", 0.1007746214909666),
    @("K7", 0.6608897862005279),
    @("M7", 0.4102477166304368),
    @("N7", 2.709023419696763),
    @("B8", 0.5489628354484637),
    @("C8", 0.1709701382410174),
    @("D8", 0.6036077724654945),
    @("E8", 0.2244374339084452),
    @("G8", 1.826266971429504),
    @("H8", 1.586746453241744),
    @("I8", 1.300019163559327),
    @("J8", 0.100954126555024),
    @("K8", 0.7102461976388952),
    @("M8", 0.4244209531544456),
    @("N8", 2.684055872100856),
    @("B9", 0.6323452151850404),
    @("C9", 0.1901261289135903),
    @("D9", 0.6176962696306418),
    @("E9", 0.2281570777805619),
    @("G9", 1.841876694645777),
    @("H9", 1.581847813760021),
    @("I9", 1.297083946432132),
    @("J9", 0.1016295415679132),
    @("K9", 0.8100838149889569),
    @("M9", 0.4540294458696224),
    @("N9", 2.640915860416008),
    @("B10", 0.6949970942623622),
    @("C10", 0.2046840274214787),
    @("D10", 0.6294002854225198),
    @("E10", 0.2313709314666355),
    @("G10", 1.857635192880053),
    @("H10", 1.581155898359185),
    @("I10", 1.297533289970211),
    @("J10", 0.1023228685515178),
    @("K10", 0.8853008053823146),
    @("M10", 0.4768909024598855),
    @("N10", 2.61277163931819),
    @("B11", 0.7238018545711213),
    @("C11", 0.2114130137195787),
    @("D11", 0.6350186440384675),
    @("E11", 0.232937467784879),
    @("G11", 1.865739334207376),
    @("H11", 1.58147236792712),
    @("I11", 1.298304529494018),
    @("J11", 0.1026810828952023),
    @("K11", 0.9199262949549336),
    @("M11", 0.4875323615314144),
    @("N11", 2.600739592646647),
    @("B12", 0.7347531055708885),
    @("C12", 0.2139764792637493),
    @("D12", 0.637188434106605),
    @("E12", 0.2335456967958436),
    @("G12", 1.868942916162126),
    @("H12", 1.581682952175299),
    @("I12", 1.298678145409099),
    @("J12", 0.1028228833100329),
    @("K12", 0.9330968481421849),
    @("M12", 0.4915967417263971),
    @("N12", 2.596294301034142),
    @("B13", 0.7323926258861775),
    @("C13", 0.2134237082166237),
    @("D13", 0.6367192532859463),
    @("E13", 0.2334140362311601),
    @("G13", 1.86824697268068),
    @("H13", 1.581633563598984),
    @("I13", 1.298594052054263),
    @("J13", 0.1027920704935354),
    @("K13", 0.9302577288497389),
    @("M13", 0.4907198628164195),
    @("N13", 2.597246736010575),
    @("B14", 0.7247019497434053),
    @("C14", 0.2116236038936279),
    @("D14", 0.6351963077566438),
    @("E14", 0.2329872063091969),
    @("G14", 1.866000193729406),
    @("H14", 1.581487874310483),
    @("I14", 1.298333632111962),
    @("J14", 0.1026926256221472),
    @("K14", 0.9210086701788782),
    @("M14", 0.4878660456304544),
    @("N14", 2.600371650807865),
    @("B15", 0.7199968488710908),
    @("C15", 0.210522987113734),
    @("D15", 0.6342689587887662),
    @("E15", 0.2327277156546188),
    @("G15", 1.864641526973486),
    @("H15", 1.581410452260542),
    @("I15", 1.298184741102816),
    @("J15", 0.1026325138616571),
    @("K15", 0.9153509853182413),
    @("M15", 0.4861225161497416),
    @("N15", 2.602300208545039),
    @("B16", 0.6931207362677867),
    @("C16", 0.2042464181321577),
    @("D16", 0.6290390264193206),
    @("E16", 0.2312706567686575),
    @("G16", 1.857124409670149),
    @("H16", 1.581147915976487),
    @("I16", 1.297494298370296),
    @("J16", 0.1023003195865044),
    @("K16", 0.883046161971663),
    @("M16", 0.4762003158249328),
    @("N16", 2.613573485271864),
    @("B17", 0.6767108459025906),
    @("C17", 0.2004232484943032),
    @("D17", 0.6259059279556709),
    @("E17", 0.2304035610299238),
    @("G17", 1.852752647328685),
    @("H17", 1.581148511664651),
    @("I17", 1.297215951426466),
    @("J17", 0.1021074921959908),
    @("K17", 0.8633328088851044),
    @("M17", 0.4701752205503098),
    @("N17", 2.620686811925175),
    @("B18", 0.667300947761106),
    @("C18", 0.1982342874816823),
    @("D18", 0.624131542129561),
    @("E18", 0.2299146711009001),
    @("G18", 1.850326169240532),
    @("H18", 1.581208269553287),
    @("I18", 1.297109203034609),
    @("J18", 0.1020006133263749),
    @("K18", 0.852032707415816),
    @("M18", 0.4667324956056618),
    @("N18", 2.624850762110277),
    @("B19", 0.6641198429579447),
    @("C19", 0.1974948639174272),
    @("D19", 0.6235355230598429),
    @("E19", 0.2297508320211463),
    @("G19", 1.849519721276806),
    @("H19", 1.581238707995055),
    @("I19", 1.297082221046274),
    @("J19", 0.1019651183340891),
    @("K19", 0.8482133044174134),
    @("M19", 0.4655707573870558),
    @("N19", 2.626273061270325),
    @("B20", 0.6784547445232647),
    @("C20", 0.2008291936856494),
    @("D20", 0.6262365862954766),
    @("E20", 0.230494846508563),
    @("G20", 1.853208915279339),
    @("H20", 1.581142299180584),
    @("I20", 1.297240060075261),
    @("J20", 0.1021276019087907),
    @("K20", 0.8654273458189721),
    @("M20", 0.4708142480817656),
    @("N20", 2.619922076282727),
    @("B21", 0.7269597085281703),
    @("C21", 0.212151921510781),
    @("D21", 0.6356424878607072),
    @("E21", 0.2331121692136051),
    @("G21", 1.866656469370525),
    @("H21", 1.581528204207359),
    @("I21", 1.298407909686595),
    @("J21", 0.102721668052304),
    @("K21", 0.9237237529898721),
    @("M21", 0.4887033399035516),
    @("N21", 2.59945077448068),
    @("B22", 0.7589140366112019),
    @("C22", 0.2196414307322812),
    @("D22", 0.6420359668333901),
    @("E22", 0.2349102564136487),
    @("G22", 1.876230563344109),
    @("H22", 1.582309331362552),
    @("I22", 1.299646612430038),
    @("J22", 0.1031457836283352),
    @("K22", 0.9621654994682274),
    @("M22", 0.5005970713155961),
    @("N22", 2.586718553719791),
    @("B23", 0.7418362834689844),
    @("C23", 0.2156359449673175),
    @("D23", 0.6386011402688041),
    @("E23", 0.2339425811996989),
    @("G23", 1.871048766111983),
    @("H23", 1.581844041035936),
    @("I23", 1.298941971641653),
    @("J23", 0.1029161455166872),
    @("K23", 0.9416172076028033),
    @("M23", 0.4942306837960544),
    @("N23", 2.593454739269177),
    @("B24", 0.6776662520423429),
    @("C24", 0.2006456377394272),
    @("D24", 0.6260870119472486),
    @("E24", 0.2304535463932851),
    @("G24", 1.853002365836659),
    @("H24", 1.581144922756295),
    @("I24", 1.297228994596729),
    @("J24", 0.1021184979095082),
    @("K24", 0.8644803017220681),
    @("M24", 0.4705252777494309),
    @("N24", 2.620267581494083),
    @("B25", 0.6095439376040872),
    @("C25", 0.1848592739629566),
    @("D25", 0.6136472973453806),
    @("E25", 0.2270663333226537),
    @("G25", 1.836901608227592),
    @("H25", 1.582662453233482),
    @("I25", 1.297420675768414),
    @("J25", 0.1014122076022126),
    @("K25", 0.7827478699853998),
    @("M25", 0.4458250998053828),
    @("N25", 2.651962983843418),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}